$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set numeric value, keep existing style untouched
function Set-Num($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Helper: cell was a text placeholder ("N/A"-style); convert to a real number
# by first copying the number-format/style from a donor numeric cell, then assigning the value.
function Set-TextToNum($addr, $val, $styleSrcAddr) {
    $ws.Range($styleSrcAddr).Copy($ws.Range($addr))
    $ws.Range($addr).Value = $val
}

# Helper: cell was numeric; convert to the shared "N/A"/"***.*" text placeholder
# by copying both the value and the style from a donor placeholder cell.
function Set-NumToText($addr, $styleSrcAddr) {
    $ws.Range($styleSrcAddr).Copy($ws.Range($addr))
}

# --- Header text updates (Volume number, report week dates) ---
$a8 = $ws.Range("A8")
$a8.Characters(21,1).Text = "10"

$c9 = $ws.Range("C9")
$c9.Characters(27,9).Text = "3/3/2025"
$c9full = $c9.Value2
$secondDateStart = $c9full.IndexOf("3/2/2025") + 1
$c9.Characters($secondDateStart,8).Text = "3/9/2025"

# --- Crime-statistics grid updates (rows 14-28) ---
Set-TextToNum "L14" -100 "K14"
Set-TextToNum "M15" 0 "K14"
Set-Num "D16" 1
Set-Num "F16" 2
Set-Num "G16" 5
Set-Num "H16" -60
Set-Num "J16" 6
Set-Num "K16" -33.333333333333
Set-Num "M16" 100
Set-Num "C17" 3
Set-Num "D17" 4
Set-Num "E17" -25
Set-Num "F17" 11
Set-Num "G17" 17
Set-Num "H17" -35.294117647058
Set-Num "I17" 24
Set-Num "J17" 23
Set-Num "K17" 4.347826086956
Set-Num "L17" 84.615384615384
Set-Num "M17" 166.666666666667
Set-Num "N17" 166.666666666667
Set-TextToNum "C18" 1 "J14"
Set-Num "F18" 5
Set-Num "H18" 150
Set-Num "I18" 15
Set-Num "K18" 400
Set-Num "L18" 7.142857142857
Set-Num "M18" -25
Set-Num "N18" -71.698113207547
Set-Num "C19" 8
Set-Num "D19" 5
Set-Num "E19" 60
Set-Num "F19" 15
Set-Num "G19" 29
Set-Num "H19" -48.275862068965
Set-Num "I19" 39
Set-Num "J19" 52
Set-Num "K19" -25
Set-Num "L19" -20.408163265306
Set-Num "M19" 39.285714285714
Set-Num "N19" 62.5
Set-NumToText "D20" "C14"
Set-NumToText "E20" "E14"
Set-Num "G20" 4
Set-Num "H20" -50
Set-Num "N20" -96.124031007751
Set-Num "C21" 12
Set-Num "D21" 10
Set-Num "E21" 20
Set-Num "F21" 36
Set-Num "H21" -36.842105263157
Set-Num "I21" 88
Set-Num "J21" 97
Set-Num "K21" -9.278350515463
Set-Num "L21" -5.376344086021
Set-Num "M21" 25.714285714285
Set-Num "N21" -61.061946902654
Set-Num "C24" 5
Set-Num "D24" 12
Set-Num "E24" -58.333333333333
Set-Num "F24" 31
Set-Num "G24" 35
Set-Num "H24" -11.428571428571
Set-Num "I24" 81
Set-Num "J24" 80
Set-Num "K24" 1.25
Set-Num "L24" -17.346938775510
Set-Num "M24" -12.903225806451
Set-Num "C25" 3
Set-Num "D25" 6
Set-Num "E25" -50
Set-Num "F25" 16
Set-Num "G25" 16
Set-Num "H25" 0
Set-Num "I25" 40
Set-Num "J25" 34
Set-Num "K25" 17.647058823529
Set-Num "L25" -31.034482758620
Set-Num "C26" 4
Set-Num "D26" 2
Set-Num "E26" 100
Set-Num "F26" 21
Set-Num "G26" 14
Set-Num "H26" 50
Set-Num "I26" 47
Set-Num "J26" 25
Set-Num "K26" 88
Set-Num "L26" 23.684210526315
Set-Num "M26" 30.555555555555
Set-NumToText "D27" "C14"
Set-NumToText "E27" "E14"
Set-Num "G27" 1
Set-Num "H27" 0
Set-NumToText "C28" "C14"
Set-Num "F28" 1
